# Apply the LinuxForHealth re-branding + version bump edits to the
# StructureDefinition-ssi-indicator workbook.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# URL
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/ssi-indicator"

# Version
$meta.Range("B3").Value = "8.0.0"

# Date
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet ---------------------------------------------------
# Clear the stray "Constraint(s)" text that had been duplicated onto the
# top-level "Extension" row (row 2) - it only belongs on the
# "Extension.extension" child row.
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""
